$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(4)
$tbl = $sh.Table
$row = $tbl.Rows.Item(1)
$cell = $row.Cells.Item(1)
$fc = $cell.Shape.Fill.ForeColor
$fc.ObjectThemeColor = 6
Write-Output ("ObjectThemeColor1: " + $fc.ObjectThemeColor)
